$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 54) continuing the existing forecast series pattern.
$ws.Range("A53:E53").Copy()
$ws.Range("A54:E54").PasteSpecial(-4122)

$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 1.049317648994741
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 0.3243937446859801
